$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.4003700119637017,
    -0.0099999997835489296,
    -0.0089999997797303166,
    -0.011999999936474381,
    -0.0059999997841035935,
    -0.0059999997781794434,
    -0.019999999732416285,
    0.057979173367773917,
    -0.0059999997738362509,
    -0.0059999997702888663,
    -0.0044999997748469411,
    -0.0059999997690303175,
    -0.0059999997653159554,
    -0.011999999744618073,
    -0.0059999997632642632,
    -0.0059999997624382573,
    -0.0059999997613608969,
    -0.008999999751322818,
    -0.077331130400462023,
    -0.0089999997840504165,
    -0.0089999997837475476,
    -0.0089999997835032985,
    -0.0089999997771048612,
    -0.067781415958468649,
    -0.041999999663034693,
    -0.0059999997772592906,
    -0.0059999997758977131,
    -0.0059999997708386488,
    -0.01199999974781818,
    -0.019999999720017314,
    -0.014999999733829483,
    -0.020999999713888329,
    -0.0059999997625501678
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
